# NIT-9012688868.xlsx — "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The account-statement table (rows 16-31) is replaced: a new worker
# (JORGE LUIS MARTINEZ QUINTANA) is inserted at the top, the existing
# OSCAR DAVID PUERTA MANJARREZ rows have their period/value columns
# reshuffled, and a new worker (ENOC DE JESUS OROZCO ARROYO, two periods)
# is appended at the bottom, ahead of the signature block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room ----------------------------------------------------
# One new row on top of the table (becomes row 16) ...
$ws.Rows(16).Insert()
# ... and two new rows at the bottom of the table (become rows 33-34),
# pushing the old footer (rows 36/37) down to 39/40.
$ws.Rows("33:34").Insert()

# --- 2. Re-apply correct formatting to the newly inserted rows -------
# Insert() copies the format of the row above, which is wrong for the
# first new row (it grabs the header style) and for the two trailing
# rows (they grab the prior last-data-row's bold/double-border style).
# Fix them up by pasting formats from a data row that already has the
# right look.
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)

$ws.Range("B17:J17").Copy()
$ws.Range("B33:J33").PasteSpecial(-4122)

# Row 34 is now the very last data row, so it should carry the bold /
# double-bottom-border style that used to belong to row 31 (WILTON's
# single row) before the table grew. Grab that look from row 32 (the
# still-bold old last row) before we normalize row 32 below.
$ws.Range("B32:J32").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122)

# Row 32 (WILTON MACIAS SANCHEZ) is no longer the last row of the
# table, so it loses the bold/double-border treatment and becomes a
# regular data row.
$ws.Range("B17:J17").Copy()
$ws.Range("B32:J32").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 3. Write the final table contents --------------------------------
$table = @(
    @(16, "CC", "73211190",   "JORGE LUIS MARTINEZ QUINTANA", "2504", 56940, 1423500),
    @(17, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2303", 43200, 1200000),
    @(18, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2302", 48000, 1200000),
    @(19, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2301", 48000, 1200000),
    @(20, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2212", 48000, 1200000),
    @(21, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2211", 48000, 1200000),
    @(22, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2210", 48000, 1200000),
    @(23, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2209", 48000, 1200000),
    @(24, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2208", 48000, 1200000),
    @(25, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2207", 48000, 1200000),
    @(26, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2206", 48000, 1200000),
    @(27, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2205", 48000, 1200000),
    @(28, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2204", 48000, 1200000),
    @(29, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2203", 48000, 1200000),
    @(30, "CC", "1043961077", "OSCAR DAVID PUERTA MANJARREZ", "2202", 48000, 1200000),
    @(31, "CE", "20445319",   "DARLENI ROSALY PEÃ?A BOLIVAR", "2307", 20107, 1160000),
    @(32, "CC", "1005259654", "WILTON MACIAS SANCHEZ",        "2307", 17013, 1160000),
    @(33, "CC", "1043005117", "ENOC DE JESUS OROZCO ARROYO",  "2504", 56940, 1423500),
    @(34, "CC", "1043005117", "ENOC DE JESUS OROZCO ARROYO",  "2503", 56940, 1423500)
)

foreach ($row in $table) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

# --- 4. Update the summary fields above the table ---------------------
$ws.Range("E11").Value = 875140   # VALOR MORA total
$ws.Range("C13").Value = 5        # Cant. Trabajadores
$ws.Range("F13").Value = 17       # Cant. Periodos
